# Applies the updated crypto price/volume figures (columns D and E, rows 2-51)
# Numeric-looking Price values are prefixed with a literal apostrophe so Excel
# keeps them as text (matching the original "d.ddd.dd"-style string formatting)
# instead of silently re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.590.55"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.468.12"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.82%  "
$ws.Range("D5").Value = "'314.60"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'91.24"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "'0.512"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").Value = "'32.57"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "'0.0793"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "2.847.71"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'6.87"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'15.82"
$ws.Range("E15").Value = "  +4.15%  "
$ws.Range("D16").Value = "2.473.57"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "'0.777"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "41.609.03"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'6.50"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "0.0₃0942"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").Value = "'71.23"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").Value = "'238.46"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'24.65"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'35.30"
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("D31").Value = "'156.01"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "'0.0758"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'17.25"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  -8.87%  "
$ws.Range("D37").Value = "'2.88"
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("D38").Value = "'0.114"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("E39").Value = "  +3.67%  "
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "1.960.29"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").Value = "'18.56"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("D46").Value = "'2.91"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  +5.35%  "
$ws.Range("D48").Value = "2.705.10"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'97.06"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "'67.34"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -1.50%  "
